$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed the "CasesTab" query: removed the erroneous `co:cohort` match/variable
# and the trailing `Cohort` output column (variable/query error fix from the
# commit message), leaving the query ending at `Response to Treatment`.
$newCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
    "WHERE demo.breed IN ['Cocker Spaniel']`n" +
    "MATCH (c)<--(diag:diagnosis)`n" +
    "OPTIONAL MATCH (samp:sample)-->(c)`n" +
    "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
    "WITH DISTINCT c, s, demo, diag, co`n" +
    "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
    "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
    "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
    "        coalesce(demo.breed, '') AS Breed ,`n" +
    "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
    "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
    "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
    "        coalesce(demo.sex, '') AS Sex ,`n" +
    "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
    "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
    "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value2 = $newCasesQuery

# Update the selected/active cell to B2 (previously B4, scrolled to A4).
$ws.Range("B2").Select()
